$d = $word.ActiveDocument

$pairs = @(
    @("44÷5=8, 4", "36÷8=4, 4"),
    @("16÷8=2, 0", "49÷3=16, 1"),
    @("77÷3=25, 2", "45÷4=11, 1"),
    @("86÷8=10, 6", "20÷9=2, 2"),
    @("24÷4=6, 0", "38÷2=19, 0"),
    @("41÷9=4, 5", "17÷4=4, 1"),
    @("14÷5=2, 4", "67÷4=16, 3"),
    @("64÷5=12, 4", "12÷7=1, 5"),
    @("57÷3=19, 0", "89÷9=9, 8"),
    @("40÷3=13, 1", "34÷2=17, 0"),
    @("11÷4=2, 3", "25÷2=12, 1"),
    @("93÷9=10, 3", "29÷9=3, 2"),
    @("34÷8=4, 2", "69÷7=9, 6"),
    @("50÷8=6, 2", "34÷7=4, 6"),
    @("71÷4=17, 3", "92÷6=15, 2"),
    @("68÷2=34, 0", "93÷9=10, 3"),
    @("38÷4=9, 2", "76÷2=38, 0"),
    @("33÷8=4, 1", "35÷7=5, 0"),
    @("60÷4=15, 0", "45÷3=15, 0"),
    @("83÷4=20, 3", "48÷6=8, 0"),
    @("53÷4=13, 1", "92÷8=11, 4"),
    @("48÷4=12, 0", "59÷2=29, 1"),
    @("38÷8=4, 6", "59÷9=6, 5"),
    @("37÷9=4, 1", "47÷5=9, 2"),
    @("73÷3=24, 1", "91÷7=13, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
